# =========================================================
# Edit script for cs-en-us-116pct.xlsx weekly crime data update
# =========================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- 1. Update the two rich-text header cells (new volume/issue number, new date range) ---
$ws.Range("A8").Value = "Volume 32   Number  41"
$ws.Range("C9").Value = "Report Covering the Week  10/6/2025  Through  10/12/2025"

# --- 2. Plain numeric value updates (counts and % changes that stay numeric) ---
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 26
$ws.Range("K15").Value = 160
$ws.Range("L15").Value = 160
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 55
$ws.Range("K16").Value = -15.384615384615
$ws.Range("L16").Value = -17.910447761194
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -18.75
$ws.Range("I17").Value = 202
$ws.Range("J17").Value = 183
$ws.Range("K17").Value = 10.382513661202
$ws.Range("L17").Value = 6.878306878306
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 0
$ws.Range("L18").Value = -29.032258064516
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -4.347826086956
$ws.Range("I19").Value = 177
$ws.Range("J19").Value = 212
$ws.Range("K19").Value = -16.509433962264
$ws.Range("L19").Value = -28.629032258064
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 130
$ws.Range("J20").Value = 139
$ws.Range("K20").Value = -6.474820143884
$ws.Range("L20").Value = -21.686746987951
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 60
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 634
$ws.Range("J21").Value = 669
$ws.Range("K21").Value = -5.231689088191
$ws.Range("L21").Value = -14.784946236559
$ws.Range("C24").Value = 10
$ws.Range("E24").Value = 11.111111111111
$ws.Range("F24").Value = 36
$ws.Range("G24").Value = 35
$ws.Range("H24").Value = 2.857142857142
$ws.Range("I24").Value = 456
$ws.Range("J24").Value = 435
$ws.Range("K24").Value = 4.827586206896
$ws.Range("L24").Value = -4.201680672268
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = -16.666666666666
$ws.Range("I25").Value = 83
$ws.Range("J25").Value = 124
$ws.Range("K25").Value = -33.064516129032
$ws.Range("L25").Value = -36.641221374045
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -38.461538461538
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = -11.764705882352
$ws.Range("I26").Value = 361
$ws.Range("J26").Value = 352
$ws.Range("K26").Value = 2.556818181818
$ws.Range("L26").Value = 17.973856209150
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 32
$ws.Range("K27").Value = 68.421052631578
$ws.Range("L27").Value = 88.235294117647
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 15
$ws.Range("K28").Value = -21.052631578947
$ws.Range("L28").Value = -34.782608695652
$ws.Range("G29").Value = 1
$ws.Range("G30").Value = 1

# --- 3. Cells that become text placeholders ("0" or "***.*") ---
#     Set the value (apostrophe-prefixed to force text), then copy the
#     number format/style from a cell that already uses that text style.
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'***.*"
$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "'***.*"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "'***.*"
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null

# --- 4. Cells that become plain numbers again (were text placeholders) ---
$ws.Range("C28").Value = 1
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
